$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.350.60"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.714.50"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.65"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5281"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06670"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2643"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.75"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07736"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.455"
$ws.Range("D13").Value = "1.949.69"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").Value = "1.717.88"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5789"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "0.0₅8171"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.66"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "27.339.68"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.59"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.009"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.641"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.40"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.028"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.712"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1205"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.223"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.15"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05365"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.294"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.481"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.381"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.634"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.848"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9515"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.400"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5870"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").Value = "1.156.02"
$ws.Range("E39").Value = "  +10.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01647"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.816"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8399"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.11"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "1.856.38"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.39"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4549"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.159"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05245"
$ws.Range("E51").Value = "  -0.17%  "
